$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$data = @(
    @(19330051920046, "COSCAHUA", "TZOYONTLE", "ALMA LIZETH", "FÍSICA I", "4AEV", 2),
    @(19330051920057, "GONZALEZ", "SANCHEZ", "JONATHAN", "FÍSICA I", "4AEV", 2),
    @(19330051920060, "HERNANDEZ", "CID", "ELIUTH ADELFO", "FÍSICA I", "4AEV", 2),
    @(19330051920441, "GARCIA", "ANTONIO", "ABRAHAM", "FÍSICA I", "4APV", 2),
    @(19330051920321, "DE JESUS", "DE LA CRUZ", "IGNACIO", "FÍSICA I", "4APV", 2),
    @(19330051920375, "HERNANDEZ", "ANTONIO", "MARIA GUADALUPE", "FÍSICA I", "4ARHV", 2),
    @(19330051920377, "HERNANDEZ", "FLORES", "PERLA", "FÍSICA I", "4ARHV", 2),
    @(19330051920382, "MAZAHUA", "IXMATLAHUA", "LUCERO", "FÍSICA I", "4ARHV", 2),
    @(18330051920113, "CASIANO", "HERNANDEZ", "DIANA MIRELLY", "TEMAS DE FÍSICA", "6APV", 2),
    @(18330051920122, "GUZMAN", "MORALES", "OSWALDO IVAN", "TEMAS DE FÍSICA", "6APV", 2),
    @(18330051920128, "MENDEZ", "SARMIENTO", "YULENNY", "TEMAS DE FÍSICA", "6APV", 2),
    @(18330051920130, "ROMERO", "PEREZ", "LUIS ENRIQUE", "TEMAS DE FÍSICA", "6APV", 2),
    @(18330051920135, "VILLA", "MACARIO", "URI YAEL", "TEMAS DE FÍSICA", "6APV", 2),
    @(18330051920044, "ANTONIO", "GAMBINO", "DAVID", "TEMAS DE FÍSICA", "6BEV", 2),
    @(19330051920056, "GONZALEZ", "DE LOS SANTOS", "MARTIN", "FÍSICA I", "4AEV", 1),
    @(19330051920307, "APALE", "ZEPAHUA", "ISRAEL", "FÍSICA I", "4APV", 1),
    @(19330051920312, "CARDENAS", "AMADOR", "KEVIN HONAM", "FÍSICA I", "4APV", 1),
    @(19330051920337, "MONTIEL", "ROJAS", "VICTOR YAHIR", "FÍSICA I", "4APV", 1),
    @(19330051920373, "GARCIA", "DIAZ", "GISELA", "FÍSICA I", "4ARHV", 1),
    @(19330051920374, "HERNANDEZ", "AGUILA", "JESUS", "FÍSICA I", "4ARHV", 1),
    @(19330051920381, "LOPEZ", "GARCIA", "MARIAM ABRIL", "FÍSICA I", "4ARHV", 1),
    @(18330051920116, "ESPIRITU", "BUSTAMANTE", "JOSE ANTONIO", "TEMAS DE FÍSICA", "6APV", 1),
    @(18330051920120, "GIL", "MARTINEZ", "ARELY", "TEMAS DE FÍSICA", "6APV", 1),
    @(18330051920121, "GONZALEZ", "REYES", "ANGEL RODRIGO", "TEMAS DE FÍSICA", "6APV", 1),
    @(18330051920125, "LIBONATTI", "FIGUEROA", "EBERTH JUVIEL", "TEMAS DE FÍSICA", "6APV", 1),
    @(17330051920466, "MARROQUIN", "ELIAS", "JORGE ANTONIO", "TEMAS DE FÍSICA", "6APV", 1),
    @(18330051920455, "PAZ", "MONTERROSAS", "DANIEL", "TEMAS DE FÍSICA", "6APV", 1)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}

Write-Host "Done writing rows"